# Insert a new weekly price record as row 198 in the "Piña" price sheet.
# Excel shifts the existing rows 198:300 down to 199:301, and we then
# populate the new row 198 with the new observation while copying the
# constant/static columns (A,B,C,E,F,G,H,I,J,K,R) from the row that used
# to occupy position 198 (now at 199) so they keep matching the rest of
# the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 198; this pushes rows
# 198-300 down to 199-301 and Excel grows the sheet's used range / the
# dimension to A1:T301 automatically.
$ws.Rows("198:198").Insert()

# Copy the columns that stay constant for every record in this sheet
# from the row just below (the row that used to be 198).
$ws.Range("A198").Value = $ws.Range("A199").Value2
$ws.Range("B198").Value = $ws.Range("B199").Value2
$ws.Range("C198").Value = $ws.Range("C199").Value2
$ws.Range("E198").Value = $ws.Range("E199").Value2
$ws.Range("F198").Value = $ws.Range("F199").Value2
$ws.Range("G198").Value = $ws.Range("G199").Value2
$ws.Range("H198").Value = $ws.Range("H199").Value2
$ws.Range("I198").Value = $ws.Range("I199").Value2
$ws.Range("J198").Value = $ws.Range("J199").Value2
$ws.Range("K198").Value = $ws.Range("K199").Value2
$ws.Range("R198").Value = $ws.Range("R199").Value2

# New data for the inserted record.
$ws.Range("D198").Value = 44813
$ws.Range("L198").Value = "Primera"
$ws.Range("M198").Value = 200
$ws.Range("N198").Value = 23000
$ws.Range("O198").Value = 23500
$ws.Range("P198").Value = 23250
$ws.Range("Q198").Value = "$/caja 12 unidades"
$ws.Range("S198").Value = 1938
$ws.Range("T198").Value = 12
